$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(42, 1).Value = "S6"
$ws.Cells.Item(42, 2).Value = "Yelena"
$ws.Cells.Item(42, 3).Value = "18/7/2013"
$ws.Cells.Item(42, 4).Value = "2013-07-18-yelena"
$ws.Cells.Item(42, 5).Value = "2013-07-18-15-11-12"
$ws.Cells.Item(42, 6).Value = "ssvep-15Hz"
$ws.Cells.Item(42, 7).Value = 15
$ws.Cells.Item(42, 8).Value = 0

$ws.Cells.Item(43, 1).Value = "S6"
$ws.Cells.Item(43, 2).Value = "Yelena"
$ws.Cells.Item(43, 3).Value = "18/7/2013"
$ws.Cells.Item(43, 4).Value = "2013-07-18-yelena"
$ws.Cells.Item(43, 5).Value = "2013-07-18-15-18-53"
$ws.Cells.Item(43, 6).Value = "hybrid-8-57Hz"
$ws.Cells.Item(43, 7).Value = 8.57
$ws.Cells.Item(43, 8).Value = 1

$ws.Cells.Item(44, 1).Value = "S6"
$ws.Cells.Item(44, 2).Value = "Yelena"
$ws.Cells.Item(44, 3).Value = "18/7/2013"
$ws.Cells.Item(44, 4).Value = "2013-07-18-yelena"
$ws.Cells.Item(44, 5).Value = "2013-07-18-15-26-03"
$ws.Cells.Item(44, 6).Value = "hybrid-10Hz"
$ws.Cells.Item(44, 7).Value = 10
$ws.Cells.Item(44, 8).Value = 1

$ws.Cells.Item(45, 1).Value = "S6"
$ws.Cells.Item(45, 2).Value = "Yelena"
$ws.Cells.Item(45, 3).Value = "18/7/2013"
$ws.Cells.Item(45, 4).Value = "2013-07-18-yelena"
$ws.Cells.Item(45, 5).Value = "2013-07-18-15-32-16"
$ws.Cells.Item(45, 6).Value = "ssvep-10Hz"
$ws.Cells.Item(45, 7).Value = 10
$ws.Cells.Item(45, 8).Value = 0

$ws.Cells.Item(46, 1).Value = "S6"
$ws.Cells.Item(46, 2).Value = "Yelena"
$ws.Cells.Item(46, 3).Value = "18/7/2013"
$ws.Cells.Item(46, 4).Value = "2013-07-18-yelena"
$ws.Cells.Item(46, 5).Value = "2013-07-18-15-43-31"
$ws.Cells.Item(46, 6).Value = "ssvep-12Hz"
$ws.Cells.Item(46, 7).Value = 12
$ws.Cells.Item(46, 8).Value = 0

$ws.Cells.Item(47, 1).Value = "S6"
$ws.Cells.Item(47, 2).Value = "Yelena"
$ws.Cells.Item(47, 3).Value = "18/7/2013"
$ws.Cells.Item(47, 4).Value = "2013-07-18-yelena"
$ws.Cells.Item(47, 5).Value = "2013-07-18-15-51-54"
$ws.Cells.Item(47, 6).Value = "hybrid-15Hz"
$ws.Cells.Item(47, 7).Value = 15
$ws.Cells.Item(47, 8).Value = 1

$ws.Cells.Item(48, 1).Value = "S6"
$ws.Cells.Item(48, 2).Value = "Yelena"
$ws.Cells.Item(48, 3).Value = "18/7/2013"
$ws.Cells.Item(48, 4).Value = "2013-07-18-yelena"
$ws.Cells.Item(48, 5).Value = "2013-07-18-15-58-42"
$ws.Cells.Item(48, 6).Value = "ssvep-8-57Hz"
$ws.Cells.Item(48, 7).Value = 8.57
$ws.Cells.Item(48, 8).Value = 0

$ws.Cells.Item(49, 1).Value = "S6"
$ws.Cells.Item(49, 2).Value = "Yelena"
$ws.Cells.Item(49, 3).Value = "18/7/2013"
$ws.Cells.Item(49, 4).Value = "2013-07-18-yelena"
$ws.Cells.Item(49, 5).Value = "2013-07-18-16-05-13"
$ws.Cells.Item(49, 6).Value = "hybrid-12Hz"
$ws.Cells.Item(49, 7).Value = 12
$ws.Cells.Item(49, 8).Value = 1

$ws.Range("L47").Select()
